$d = $word.ActiveDocument

# Locate the paragraph that ends the bibliography entry (the "...3040P."
# paragraph). The three paragraphs that follow it -- a blank paragraph,
# "Ver no Jupiter Salvar em pdf Salvar em docx", and the "© 2020 ..."
# footer line -- are the ones being removed by this edit.
$marker = $d.Content.Find
$marker.Text = "Risk management for hazardous chemicals. CRC Press USA:1997. 3040P."
$marker.Execute() | Out-Null

$anchorPara = $d.Range($marker.Parent.Start, $marker.Parent.Start).Paragraphs(1)
$startPara = $anchorPara.Next()          # blank paragraph right after the bibliography text
$endPara = $startPara.Next().Next()      # the "© 2020 ..." paragraph

$doomed = $d.Range($startPara.Range.Start, $endPara.Range.End)
$doomed.Delete()
